$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2383824.5
$ws.Range("I15").Value = 2383824.5
$ws.Range("K15").Value = 7151473.5
$ws.Range("M15").Value = -7151304.5
$ws.Range("H74").Value = 4028.111
$ws.Range("I74").Value = 4266.3335
$ws.Range("J74").Value = 3909
$ws.Range("K74").Value = 4266.3335
$ws.Range("L74").Value = 3909
$ws.Range("M74").Value = -3330.3335
$ws.Range("N74").Value = -5781
$ws.Range("H77").Value = 4028.111
$ws.Range("I77").Value = 4266.3335
$ws.Range("J77").Value = 3909
$ws.Range("K77").Value = 21331.6675
$ws.Range("L77").Value = 19545
$ws.Range("M77").Value = -16651.6675
$ws.Range("N77").Value = -28905
$ws.Range("H87").Value = 24828.572
$ws.Range("J87").Value = 24828.572
$ws.Range("L87").Value = 24828.572
$ws.Range("N87").Value = -27324.572
$ws.Range("H90").Value = 24828.572
$ws.Range("J90").Value = 24828.572
$ws.Range("L90").Value = 74485.716
$ws.Range("N90").Value = -86965.716
$ws.Range("H109").Value = 39092
$ws.Range("J109").Value = 39092
$ws.Range("L109").Value = 39092
$ws.Range("N109").Value = -41866
$ws.Range("H129").Value = 977.7632
$ws.Range("J129").Value = 1027.2858
$ws.Range("L129").Value = 3081.8574
$ws.Range("N129").Value = -13081.8574
$ws.Range("H132").Value = 3795.3044
$ws.Range("I132").Value = 3704.8108
$ws.Range("J132").Value = 4167.3335
$ws.Range("K132").Value = 11114.4324
$ws.Range("L132").Value = 12502.0005
$ws.Range("M132").Value = -8584.432400000002
$ws.Range("N132").Value = -17562.0005

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 34196.668
$ws.Range("J80").Value = 35971.25
$ws.Range("L80").Value = 35971.25
$ws.Range("N80").Value = -37967.25
$ws.Range("H83").Value = 34196.668
$ws.Range("J83").Value = 35971.25
$ws.Range("L83").Value = 107913.75
$ws.Range("N83").Value = -117897.75
$ws.Range("H110").Value = 1062.3846
$ws.Range("I110").Value = 856.8333
$ws.Range("J110").Value = 1238.5714
$ws.Range("K110").Value = 856.8333
$ws.Range("L110").Value = 1238.5714
$ws.Range("M110").Value = 1188.1667
$ws.Range("N110").Value = -5328.5714
$ws.Range("H122").Value = 1645.9166
$ws.Range("I122").Value = 1166.9412
$ws.Range("J122").Value = 2809.1428
$ws.Range("K122").Value = 3500.8236
$ws.Range("L122").Value = 8427.428400000001
$ws.Range("M122").Value = -1050.8236
$ws.Range("N122").Value = -13327.4284

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1793.7142
$ws.Range("I20").Value = 1854.7693
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 1854.7693
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = -1607.7693
$ws.Range("N20").Value = -1494
$ws.Range("H134").Value = 4123.081
$ws.Range("I134").Value = 4367.9375
$ws.Range("K134").Value = 13103.8125
$ws.Range("M134").Value = -10568.8125

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 724.8889
$ws.Range("I16").Value = 729.8
$ws.Range("J16").Value = 718.75
$ws.Range("K16").Value = 729.8
$ws.Range("L16").Value = 718.75
$ws.Range("M16").Value = -442.8
$ws.Range("N16").Value = -1292.75
$ws.Range("H113").Value = 724.8889
$ws.Range("I113").Value = 729.8
$ws.Range("J113").Value = 718.75
$ws.Range("K113").Value = 729.8
$ws.Range("L113").Value = 718.75
$ws.Range("M113").Value = 1440.2
$ws.Range("N113").Value = -5058.75
$ws.Range("H122").Value = 2927.5
$ws.Range("I122").Value = 2831.818
$ws.Range("K122").Value = 8495.454000000002
$ws.Range("M122").Value = -6045.454000000002
$ws.Range("H140").Value = 39272
$ws.Range("J140").Value = 39272
$ws.Range("L140").Value = 39272
$ws.Range("N140").Value = -49632

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1093.4736
$ws.Range("I5").Value = 663.3333
$ws.Range("J5").Value = 1292
$ws.Range("K5").Value = 1989.9999
$ws.Range("L5").Value = 3876
$ws.Range("M5").Value = -1877.9999
$ws.Range("N5").Value = -4100
$ws.Range("H12").Value = 54.333332
$ws.Range("I12").Value = 61.88889
$ws.Range("J12").Value = 50.555557
$ws.Range("K12").Value = 185.66667
$ws.Range("L12").Value = 151.666671
$ws.Range("M12").Value = -12.66667000000001
$ws.Range("N12").Value = -497.666671
$ws.Range("H113").Value = 1545.8572
$ws.Range("I113").Value = 2082.4443
$ws.Range("J113").Value = 580
$ws.Range("K113").Value = 6247.3329
$ws.Range("L113").Value = 1740
$ws.Range("M113").Value = -4077.3329
$ws.Range("N113").Value = -6080
$ws.Range("H135").Value = 1093.4736
$ws.Range("I135").Value = 663.3333
$ws.Range("J135").Value = 1292
$ws.Range("K135").Value = 5969.9997
$ws.Range("L135").Value = 11628
$ws.Range("M135").Value = -3434.9997
$ws.Range("N135").Value = -16698

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2993.3333
$ws.Range("I80").Value = 2980
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2980
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1982
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2993.3333
$ws.Range("I83").Value = 2980
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 14900
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -9908
$ws.Range("N83").Value = -24984
$ws.Range("H122").Value = 1999.9524
$ws.Range("I122").Value = 1816.2307
$ws.Range("J122").Value = 2298.5
$ws.Range("K122").Value = 5448.6921
$ws.Range("L122").Value = 6895.5
$ws.Range("M122").Value = -2998.6921
$ws.Range("N122").Value = -11795.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 481
$ws.Range("I46").Value = 483
$ws.Range("K46").Value = 483
$ws.Range("M46").Value = -295
$ws.Range("H61").Value = 1568.1818
$ws.Range("I61").Value = 1445
$ws.Range("K61").Value = 1445
$ws.Range("M61").Value = -1243
$ws.Range("H82").Value = 2526.818
$ws.Range("I82").Value = 1400
$ws.Range("J82").Value = 2777.2222
$ws.Range("K82").Value = 1400
$ws.Range("L82").Value = 2777.2222
$ws.Range("M82").Value = -1039
$ws.Range("N82").Value = -3499.2222
$ws.Range("H85").Value = 2526.818
$ws.Range("I85").Value = 1400
$ws.Range("J85").Value = 2777.2222
$ws.Range("K85").Value = 1400
$ws.Range("L85").Value = 2777.2222
$ws.Range("M85").Value = -152
$ws.Range("N85").Value = -5273.2222
$ws.Range("H113").Value = 1568.1818
$ws.Range("I113").Value = 1445
$ws.Range("K113").Value = 1445
$ws.Range("M113").Value = 725
$ws.Range("H132").Value = 2365.595
$ws.Range("I132").Value = 1932.5834
$ws.Range("J132").Value = 3733
$ws.Range("K132").Value = 5797.7502
$ws.Range("L132").Value = 11199
$ws.Range("M132").Value = -3267.7502
$ws.Range("N132").Value = -16259
